# Refresh the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# with updated snapshot values, mirroring the scheduled GitHub Actions update.
#
# Note: several "Price" cells contain text that looks numeric (e.g. "8.10",
# "355.00"). Assigning such text directly to .Value would make Excel coerce
# it to a real number and drop the significant trailing zero. To keep these
# cells as text (matching the original inlineStr cells), we prefix the value
# with a leading apostrophe, Excel's standard "treat as text" marker; Excel
# strips that marker from the stored value automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.988.77'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '2.526.29'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''589.05'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('D6').Value = '''172.53'
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '2.526.27'
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('E10').Value = '  +0.48%  '
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('D12').Value = '''5.13'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '''0.342'
$ws.Range('E13').Value = '  -3.72%  '
$ws.Range('D14').Value = '''26.52'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '2.989.64'
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '66.891.80'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = '2.524.47'
$ws.Range('E18').Value = '  -2.64%  '
$ws.Range('D19').Value = '''8.10'
$ws.Range('E19').Value = '  +4.71%  '
$ws.Range('D20').Value = '''11.34'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').Value = '''355.00'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').Value = '''4.17'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = '''4.61'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = '''1.98'
$ws.Range('E24').Value = '  +5.29%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').Value = '''69.75'
$ws.Range('D27').Value = '''9.98'
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = '2.653.81'
$ws.Range('D30').Value = '0.0₃0973'
$ws.Range('D31').Value = '''531.67'
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').Value = '''157.01'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').Value = '''18.58'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').Value = '''18.44'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('D41').Value = '''0.354'
$ws.Range('E41').Value = '  -1.73%  '
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('E45').Value = '  +3.11%  '
$ws.Range('D46').Value = '''149.09'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '''0.555'
$ws.Range('E47').Value = '  -2.09%  '
$ws.Range('D48').Value = '0.0₆0277'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').Value = '''1.68'
$ws.Range('E51').Value = '  -0.21%  '
